# Correccion de Importes Informe Sigehos
# Appends 17 new "Financiador" rows (with their "Tipo Cobertura") to the
# bottom of the Hoja1 listing, right after the existing last row (1022).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# $null placeholder marks the one row (1032) whose "Financiador" text
# genuinely starts with a literal apostrophe; it is written separately,
# in-place/in-order, further below so the shared-string table keeps the
# same left-to-right allocation order as the rows themselves.
$financiadores = @(
    'CONEXION SALUD - CONEXION SALUD SRL',
    'OSPIC SANTA FE - OBRA SOCIAL DEL PERSONAL DE LA INDUSTRIA DEL CAUCHO DE SANTA FE',
    'SANATORIO GREYT - SANATORIO GREYTON S.A.',
    'INCLUIR TIERRA - INCLUIR SALUD TIERRA DEL FUEGO ANTARTIDA E ISLAS DEL ATLANTICO SUR',
    'GERMED S.A. (LA - GERMED S.A. (LA PEQUEÑA FAMILIA MEDICINA PREPAGA)',
    'IOSE - INSTITUTO OBRA SOCIAL DEL EJERCITO',
    'EMBRACE SRL - EMBRACE SRL',
    'APSOT - OBRA SOCIAL ASOCIACION DEL PERSONAL SUPERIOR DE LA ORGANIZACION TECHINT',
    'HTAL BRITANICO - Plan de Salud Hospital Britanico de Bs. As.',
    $null,
    'OSPA VIAL - OBRA SOCIAL DEL PERSONAL DE LA ACTIVIDAD VIAL',
    'EQUITATIVA DEL - LA EQUITATIVA DEL PLATA SA DE SEGUROS',
    'AGROSALTA COOP. - AGROSALTA COOP. SEGUROS LIMITADA',
    'BEST DOCTORS IN - BEST DOCTORS INSURANCE SERVICES',
    'FESTIQYPRA - OBRA SOCIAL DEL PERSONAL DE LA FEDERACION DE SINDICATOS DE LA INDUSTRIA QUIMICA Y PETROQUIMICA DE LA REPUBLICA ARGENTINA',
    'INST MED MODELO - PRESTACIONES MEDICAS INTEGRALES - INSTITUTO MEDICO MODELO S.A.',
    '- CAJA DE PREV. SOCIAL DE PROF. DE ING. DE SANTA FE'
)

$coberturas = @(
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'Incluir Extra Cápita',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas',
    'OOSS y Prepagas'
)

$firstNewRow = 1023
$lastDataRow = $firstNewRow + $financiadores.Count - 1

for ($i = 0; $i -lt $financiadores.Count; $i++) {
    $r = $firstNewRow + $i
    # Clone the formatting of the previous row (font, fill, number format,
    # row height) down onto the new row before writing its values.
    $srcRow = $r - 1
    $ws.Range("A" + $srcRow + ":B" + $srcRow).Copy($ws.Range("A" + $r + ":B" + $r))
    if ($financiadores[$i] -ne $null) {
        $ws.Range("A$r").Value = $financiadores[$i]
    }
    $ws.Range("B$r").Value = $coberturas[$i]
}

# Row 1032's "Financiador" text genuinely starts with a literal apostrophe
# character. Assigning that straight to .Value would make Excel treat the
# leading quote as a text-qualifier marker (quote-prefix) instead of as
# part of the content, so we build it as a formula first (which preserves
# the literal character) and then convert that formula to a plain value.
$apoRow = 1032
$ws.Range("A$apoRow").Formula = '="''- CAJA DE PREV. SOCIAL DE PROF. DE ING. DE SANTA FE"'
$ws.Range("A$apoRow").Copy()
$ws.Range("A$apoRow").PasteSpecial(-4163)

# Update the sheet dimension / selection to reflect the newly added rows.
$ws.Range("A1043").Select()
